$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their text formatting so numeric-looking
# strings (e.g. "1.00", "0.996") are not silently converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '55.426.39'
$ws.Range("E2").Value = '  +1.74%  '
$ws.Range("D3").Value = '2.297.87'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '507.25'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '130.21'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '0.532'
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").Value = '2.322.10'
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").Value = '0.0985'
$ws.Range("E10").Value = '  +2.58%  '
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").Value = '5.09'
$ws.Range("E12").Value = '  +7.35%  '
$ws.Range("E13").Value = '  +0.91%  '
$ws.Range("D14").Value = '23.95'
$ws.Range("E14").Value = '  +4.62%  '
$ws.Range("D15").Value = '2.706.31'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("D16").Value = '55.258.45'
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("D18").Value = '2.307.76'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("D19").Value = '10.78'
$ws.Range("E19").Value = '  +5.14%  '
$ws.Range("D20").Value = '4.20'
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").Value = '312.12'
$ws.Range("E21").Value = '  +2.42%  '
$ws.Range("D22").Value = '6.63'
$ws.Range("E22").Value = '  +3.76%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = '60.62'
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").Value = '0.152'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '7.54'
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").Value = '172.71'
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").Value = '6.16'
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").Value = '0.0₃0712'
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("D31").Value = '1.64'
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  +4.74%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '18.11'
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '1.24'
$ws.Range("E36").Value = '  +2.68%  '
$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").Value = '0.922'
$ws.Range("E37").Value = '  -4.23%  '
$ws.Range("D38").Value = '3.91'
$ws.Range("E38").Value = '  +4.43%  '
$ws.Range("D39").Value = '36.82'
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").Value = '1.45'
$ws.Range("E40").Value = '  +2.28%  '
$ws.Range("D41").Value = '0.378'
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").Value = '135.62'
$ws.Range("E42").Value = '  +7.57%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '3.45'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '4.93'
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("D45").Value = '262.11'
$ws.Range("E45").Value = '  +7.79%  '
$ws.Range("D46").Value = '0.0506'
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '0.0914'
$ws.Range("E47").Value = '  +1.88%  '
$ws.Range("D48").Value = '0.553'
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("D49").Value = '0.378'
$ws.Range("E49").Value = '  +0.83%  '
$ws.Range("D50").Value = '0.0211'
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '16.54'
$ws.Range("E51").Value = '  +0.85%  '
